$d = $word.ActiveDocument

$replacements = @(
    @("867×3=", "297×5="),
    @("701×9=", "869×2="),
    @("932×3=", "502×6="),
    @("213×2=", "362×7="),
    @("290×3=", "160×9="),
    @("721×2=", "586×9="),
    @("459×7=", "339×2="),
    @("909×8=", "816×7="),
    @("641×2=", "591×7="),
    @("813×3=", "394×7="),
    @("823×2=", "253×3="),
    @("684×6=", "928×8="),
    @("909×6=", "571×4="),
    @("385×5=", "116×8="),
    @("391×2=", "212×3="),
    @("254×5=", "608×9="),
    @("993×8=", "269×2="),
    @("432×6=", "254×6="),
    @("758×8=", "738×6="),
    @("551×9=", "535×6="),
    @("675×4=", "644×7="),
    @("318×4=", "259×9="),
    @("795×7=", "686×6="),
    @("162×8=", "247×9="),
    @("836×5=", "767×7=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
